# Task H: Vehicle Menu
# - Confirming a vehicle purchase removes it from the Vehicles list (the
#   "20 Berrari" row is spent/cleared) and payment is taken.
# - The active sheet/tab switches to "Vehicles" with the now-empty row
#   selected.

$wb = $excel.ActiveWorkbook

$wsVehicles = $wb.Worksheets.Item("Vehicles")

# Clear out the purchased vehicle's row (A5:C5) - price cell (B5) keeps its
# currency formatting but becomes blank, name/stat cells are removed.
$wsVehicles.Range("A5").ClearContents()
$wsVehicles.Range("B5").ClearContents()
$wsVehicles.Range("C5").ClearContents()

# Switch to the Vehicles tab and select the now-empty row, matching the
# in-game camera/menu focus change described in the commit message.
$wsVehicles.Activate()
$wsVehicles.Range("A5:C5").Select()
